# Rename the original sheet and add the new "Normalisation" sheet after it.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Comparison"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Normalisation"

# Header row
$ws2.Range("A1").Value = "Action"
$ws2.Range("B1").Value = "Time"
$ws2.Range("C1").Value = "Content"

# Data rows
$ws2.Range("A2").Value = "Upload CSV"
$ws2.Range("B2").Value = "5 min"
$ws2.Range("C2").Value = "df = pd.read_csv('file.csv')"

$ws2.Range("A3").Value = "Visualize Data"
$ws2.Range("B3").Value = "10 min"
$ws2.Range("C3").Value = "df.describe()"

$ws2.Range("A4").Value = "Apply Normalization"
$ws2.Range("B4").Value = "5 min"
$ws2.Range("C4").Value = "from sklearn.preprocessing import MinMaxScaler <br> scaler = MinMaxScaler() <br> df_scaled = pd.DataFrame(scaler.fit_transform(df), columns=df.columns)"

$ws2.Range("A5").Value = "Verify Normalization"
$ws2.Range("B5").Value = "2 min"
$ws2.Range("C5").Value = "df_scaled.describe()"

# Overall row
$ws2.Range("A6").Value = "Overall"
$ws2.Range("B6").Value = "22 min"

# Formatting: bold size-13 for header/overall rows, regular size-13 for body rows.
# Seed the bold look by copying the bold (Task header) format from Comparison!A1,
# then shrink it to size 13 -- this yields the bold/13 font before the plain/13
# font, matching how the two new fonts end up ordered in the style table.
$ws1.Range("A1").Copy()
$ws2.Range("A1:C1").PasteSpecial(-4122)
$ws2.Range("A1:C1").Font.Size = 13

$ws1.Range("A1").Copy()
$ws2.Range("A6:B6").PasteSpecial(-4122)
$ws2.Range("A6:B6").Font.Size = 13

$ws2.Range("A2:C5").Font.Size = 13
